$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D36").Value = "Deep Neural Networks with Noisy Labels"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/377"

$ws.Range("D50").Value = "위상정렬 (topological sorting)"
